$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the footer/metadata rows 878-882
$ws.Rows("878:882").Delete()

# Update header row to snake_case English column names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Capitalize Spanish connector words (de/del/el/la/los/y) in place names
# and fix "MonteMorelos" -> "Montemorelos"
$ws.Range("B8").Value = "Pabellón De Arteaga"
$ws.Range("B9").Value = "Rincón De Romos"
$ws.Range("B10").Value = "San Francisco De Los Romo"
$ws.Range("B11").Value = "San José De Gracia"
$ws.Range("B31").Value = "Comitán De Domínguez"
$ws.Range("B47").Value = "Mazapa De Madero"
$ws.Range("B49").Value = "Ocozocoautla De Espinosa"
$ws.Range("B56").Value = "Salto De Agua"
$ws.Range("B89").Value = "Guadalupe Y Calvo"
$ws.Range("B91").Value = "Hidalgo Del Parral"
$ws.Range("B105").Value = "San Francisco Del Oro"
$ws.Range("B122").Value = "San Juan De Sabinas"
$ws.Range("A132").Value = "Ciudad De México"
$ws.Range("B136").Value = "Cuajimalpa De Morelos"
$ws.Range("B150").Value = "Coneto De Comonfort"
$ws.Range("B164").Value = "Nombre De Dios"
$ws.Range("B168").Value = "Pánuco De Coronado"
$ws.Range("B174").Value = "San Juan De Guadalupe"
$ws.Range("B175").Value = "San Luis Del Cordero"
$ws.Range("B176").Value = "San Pedro Del Gallo"
$ws.Range("A183").Value = "Estado De México"
$ws.Range("B183").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B187").Value = "Atizapán De Zaragoza"
$ws.Range("B196").Value = "Ecatepec De Morelos"
$ws.Range("B200").Value = "Ixtapan De La Sal"
$ws.Range("B210").Value = "Naucalpan De Juárez"
$ws.Range("B219").Value = "Tenango Del Valle"
$ws.Range("B223").Value = "Tlalnepantla De Baz"
$ws.Range("B234").Value = "Apaseo El Alto"
$ws.Range("B235").Value = "Apaseo El Grande"
$ws.Range("B242").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B246").Value = "Jaral Del Progreso"
$ws.Range("B253").Value = "Purísima Del Rincón"
$ws.Range("B257").Value = "San Diego De La Unión"
$ws.Range("B259").Value = "San Francisco Del Rincón"
$ws.Range("B260").Value = "San Luis De La Paz"
$ws.Range("B261").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B262").Value = "Silao De La Victoria"
$ws.Range("B266").Value = "Valle De Santiago"
$ws.Range("B271").Value = "Acapulco De Juárez"
$ws.Range("B273").Value = "Ajuchitlán Del Progreso"
$ws.Range("B276").Value = "Atoyac De Álvarez"
$ws.Range("B277").Value = "Ayutla De Los Libres"
$ws.Range("B278").Value = "Chilapa De Álvarez"
$ws.Range("B279").Value = "Chilpancingo De Los Bravo"
$ws.Range("B280").Value = "Coyuca De Benítez"
$ws.Range("B281").Value = "Coyuca De Catalán"
$ws.Range("B283").Value = "Cutzamala De Pinzón"
$ws.Range("B288").Value = "Iguala De La Independencia"
$ws.Range("B289").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B290").Value = "Zihuatanejo De Azueta"
$ws.Range("B291").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B301").Value = "Taxco De Alarcón"
$ws.Range("B302").Value = "Técpan De Galeana"
$ws.Range("B305").Value = "Tlalixtaquilla De Maldonado"
$ws.Range("B306").Value = "Tlapa De Comonfort"
$ws.Range("B321").Value = "Jacala De Ledezma"
$ws.Range("B323").Value = "Molango De Escamilla"
$ws.Range("B325").Value = "Nopala De Villagrán"
$ws.Range("B326").Value = "Pachuca De Soto"
$ws.Range("B328").Value = "Progreso De Obregón"
$ws.Range("B331").Value = "Tenango De Doria"
$ws.Range("B334").Value = "Tula De Allende"
$ws.Range("B335").Value = "Tulancingo De Bravo"
$ws.Range("B336").Value = "Zacualtipán De Ángeles"
$ws.Range("B345").Value = "Atotonilco El Alto"
$ws.Range("B346").Value = "Autlán De Navarro"
$ws.Range("B349").Value = "Cañadas De Obregón"
$ws.Range("B352").Value = "Concepción De Buenos Aires"
$ws.Range("B353").Value = "Cuautitlán De García Barragán"
$ws.Range("B358").Value = "Encarnación De Díaz"
$ws.Range("B361").Value = "Ixtlahuacán Del Río"
$ws.Range("B367").Value = "Lagos De Moreno"
$ws.Range("B372").Value = "Ojuelos De Jalisco"
$ws.Range("B377").Value = "San Cristóbal De La Barranca"
$ws.Range("B379").Value = "San Juan De Los Lagos"
$ws.Range("B380").Value = "San Juanito De Escobedo"
$ws.Range("B382").Value = "San Martín De Bolaños"
$ws.Range("B383").Value = "San Miguel El Alto"
$ws.Range("B386").Value = "Tamazula De Gordiano"
$ws.Range("B389").Value = "Tepatitlán De Morelos"
$ws.Range("B394").Value = "Unión De San Antonio"
$ws.Range("B395").Value = "Unión De Tula"
$ws.Range("B396").Value = "Valle De Juárez"
$ws.Range("B398").Value = "Yahualica De González Gallo"
$ws.Range("B401").Value = "Zapotlán El Grande"
$ws.Range("B415").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B457").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B473").Value = "Coatlán Del Río"
$ws.Range("B477").Value = "Puente De Ixtla"
$ws.Range("B481").Value = "Tetela Del Volcán"
$ws.Range("B482").Value = "Tlaltizapán De Zapata"
$ws.Range("B491").Value = "Ixtlán Del Río"
$ws.Range("B508").Value = "Mier Y Noriega"
$ws.Range("B509").Value = "Montemorelos"
$ws.Range("B512").Value = "San Nicolás De Los Garza"
$ws.Range("B518").Value = "Chalcatongo De Hidalgo"
$ws.Range("B521").Value = "Guevea De Humboldt"
$ws.Range("B522").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B523").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B524").Value = "Ixtlán De Juárez"
$ws.Range("B525").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B528").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B529").Value = "Oaxaca De Juárez"
$ws.Range("B530").Value = "Ocotlán De Morelos"
$ws.Range("B531").Value = "Putla Villa De Guerrero"
$ws.Range("B533").Value = "San Antonino El Alto"
$ws.Range("B534").Value = "San Dionisio Del Mar"
$ws.Range("B548").Value = "San Miguel Del Puerto"
$ws.Range("B583").Value = "Tataltepec De Valdés"
$ws.Range("B584").Value = "Villa De Etla"
$ws.Range("B585").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B586").Value = "Zimatlán De Álvarez"
$ws.Range("B592").Value = "Chalchicomula De Sesma"
$ws.Range("B606").Value = "Mazapiltepec De Juárez"
$ws.Range("B608").Value = "Palmar De Bravo"
$ws.Range("B618").Value = "Tepexi De Rodríguez"
$ws.Range("B621").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B631").Value = "Amealco De Bonfil"
$ws.Range("B633").Value = "Cadereyta De Montes"
$ws.Range("B636").Value = "Jalpan De Serra"
$ws.Range("B637").Value = "Landa De Matamoros"
$ws.Range("B639").Value = "Pinal De Amoles"
$ws.Range("B641").Value = "San Juan Del Río"
$ws.Range("B649").Value = "Axtla De Terrazas"
$ws.Range("B655").Value = "Ciudad Del Maíz"
$ws.Range("B663").Value = "Mexquitic De Carmona"
$ws.Range("B667").Value = "San Ciro De Acosta"
$ws.Range("B672").Value = "Santa María Del Río"
$ws.Range("B673").Value = "Soledad De Graciano Sánchez"
$ws.Range("B681").Value = "Villa De Arista"
$ws.Range("B682").Value = "Villa De Arriaga"
$ws.Range("B683").Value = "Villa De Guadalupe"
$ws.Range("B684").Value = "Villa De Ramos"
$ws.Range("B685").Value = "Villa De Reyes"
$ws.Range("B737").Value = "Soto La Marina"
$ws.Range("B746").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B747").Value = "San Pablo Del Monte"
$ws.Range("B758").Value = "Amatlán De Los Reyes"
$ws.Range("B761").Value = "Boca Del Río"
$ws.Range("B764").Value = "Castillo De Teayo"
$ws.Range("B765").Value = "Cazones De Herrera"
$ws.Range("B771").Value = "Cosamaloapan De Carpio"
$ws.Range("B776").Value = "Hueyapan De Ocampo"
$ws.Range("B777").Value = "Ignacio De La Llave"
$ws.Range("B779").Value = "Ixhuatlán Del Café"
$ws.Range("B780").Value = "Ixhuatlán Del Sureste"
$ws.Range("B785").Value = "Juchique De Ferrer"
$ws.Range("B789").Value = "Martínez De La Torre"
$ws.Range("B791").Value = "Medellín De Bravo"
$ws.Range("B794").Value = "Nanchital De Lázaro Cárdenas Del Río"
$ws.Range("B798").Value = "Ozuluama De Mascareñas"
$ws.Range("B802").Value = "Paso De Ovejas"
$ws.Range("B804").Value = "Poza Rica De Hidalgo"
$ws.Range("B809").Value = "Sayula De Alemán"
$ws.Range("B812").Value = "Tatahuicapan De Juárez"
$ws.Range("B836").Value = "Cañitas De Felipe Pescador"
$ws.Range("B838").Value = "Concepción Del Oro"
$ws.Range("B855").Value = "Nochistlán De Mejía"
$ws.Range("B865").Value = "Teúl De González Ortega"
$ws.Range("B866").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B870").Value = "Villa De Cos"

# Correct float rounding in percentage column (last-digit precision fix)
$ws.Range("D6").Value = 0.0009376465072667604
$ws.Range("D12").Value = 0.0009376465072667604
$ws.Range("D32").Value = 0.0009376465072667604
$ws.Range("D71").Value = 0.0009376465072667604
$ws.Range("D89").Value = 0.0009376465072667604
$ws.Range("D122").Value = 0.0009376465072667604
$ws.Range("D131").Value = 0.0009376465072667604
$ws.Range("D143").Value = 0.0009376465072667604
$ws.Range("D163").Value = 0.0009376465072667604
$ws.Range("D165").Value = 0.0009376465072667604
$ws.Range("D166").Value = 0.0009376465072667604
$ws.Range("D225").Value = 0.0009376465072667604
$ws.Range("D279").Value = 0.0009376465072667604
$ws.Range("D319").Value = 0.0009376465072667604
$ws.Range("D343").Value = 0.0009376465072667604
$ws.Range("D354").Value = 0.0009376465072667604
$ws.Range("D390").Value = 0.0009376465072667604
$ws.Range("D440").Value = 0.0009376465072667604
$ws.Range("D447").Value = 0.0009376465072667604
$ws.Range("D453").Value = 0.0009376465072667604
$ws.Range("D486").Value = 0.009220190654789812
$ws.Range("D525").Value = 0.0009376465072667604
$ws.Range("D529").Value = 0.0009376465072667604
$ws.Range("D530").Value = 0.0009376465072667604
$ws.Range("D537").Value = 0.0009376465072667604
$ws.Range("D635").Value = 0.0009376465072667604
$ws.Range("D644").Value = 0.0009376465072667604
$ws.Range("D645").Value = 0.0009376465072667604
$ws.Range("D668").Value = 0.009220190654789812
$ws.Range("D704").Value = 0.0009376465072667604
$ws.Range("D706").Value = 0.0009376465072667604
$ws.Range("D746").Value = 0.0009376465072667604
$ws.Range("D771").Value = 0.0009376465072667604
$ws.Range("D818").Value = 0.0009376465072667604
$ws.Range("D849").Value = 0.0009376465072667604

